# [Site Web V4] amélioration du header et des listes
#
# Rewrites the two most recent task-log entries on "Feuil1":
#   - row 19: text, date (42329 -> 42330) and duration (0.45 -> 3)
#   - row 20: text, date (42329 -> 42330) and duration (2 -> 0.5)
# and moves the active selection from I13 to D20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("B19").Value = "Réorganisation de page / Centrage de certains élément/ refactoring / amélioration des listes"
$ws.Range("C19").Value = 42330
$ws.Range("D19").Value = 3

$ws.Range("B20").Value = "amélioration du header et refacto avec un peu de gestion de connexion d'un compte"
$ws.Range("C20").Value = 42330
$ws.Range("D20").Value = 0.5

$ws.Range("D20").Select()
